$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking strings (e.g. "1.00", "533.90")
# keep their exact formatting instead of being parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '58.941.69'
$ws.Range("E2").Value = '  +1.23%  '
$ws.Range("D3").Value = '2.502.86'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '533.90'
$ws.Range("E5").Value = '  +2.56%  '
$ws.Range("D6").Value = '134.05'
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").Value = '0.570'
$ws.Range("E8").Value = '  +2.60%  '
$ws.Range("D9").Value = '2.506.45'
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("D10").Value = '0.0994'
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").Value = '  -2.82%  '
$ws.Range("D12").Value = '5.17'
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("D13").Value = '0.330'
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("D14").Value = '2.949.71'
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("D15").Value = '58.755.41'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").Value = '22.34'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '2.503.99'
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").Value = '10.61'
$ws.Range("E19").Value = '  -1.20%  '
$ws.Range("E20").Value = '  +1.87%  '
$ws.Range("D21").Value = '321.10'
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("D22").Value = '6.17'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").Value = '65.84'
$ws.Range("E24").Value = '  +3.75%  '
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").Value = '0.160'
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("D28").Value = '7.44'
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("D29").Value = '0.0₃0757'
$ws.Range("E29").Value = '  +1.02%  '
$ws.Range("D30").Value = '172.52'
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("D31").Value = '1.73'
$ws.Range("E31").Value = '  +1.46%  '
$ws.Range("D32").Value = '6.27'
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("D36").Value = '18.11'
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("E37").Value = '  -3.84%  '
$ws.Range("D38").Value = '3.96'
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("E39").Value = '  +3.60%  '
$ws.Range("D40").Value = '0.828'
$ws.Range("E40").Value = '  +6.24%  '
$ws.Range("D41").Value = '36.38'
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("D42").Value = '3.47'
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("D43").Value = '275.01'
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("D44").Value = '131.70'
$ws.Range("E44").Value = '  +6.53%  '
$ws.Range("D45").Value = '5.02'
$ws.Range("E45").Value = '  -2.37%  '
$ws.Range("D46").Value = '0.593'
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("E47").Value = '  +1.57%  '
$ws.Range("D48").Value = '0.0510'
$ws.Range("E48").Value = '  +2.51%  '
$ws.Range("D49").Value = '0.0217'
$ws.Range("E49").Value = '  +2.00%  '
$ws.Range("D50").Value = '16.75'
$ws.Range("E50").Value = '  -1.90%  '
$ws.Range("D51").Value = '1.753.27'
$ws.Range("E51").Value = '  +0.44%  '

# Restore the default (unstyled) look for column D so no stray number format
# is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
